# Pre-test per piu' zone coperte
# Replace rows 3-9 with the updated scheduling data (covering more zones),
# and remove the now-obsolete trailing rows 10-14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 7,14

$data[0,0] = 251346
$data[0,1] = "CASON"
$data[0,2] = 36.5
$data[0,3] = 70.16363636363636
$data[0,4] = "2025-04-10 10:41:00"
$data[0,5] = "2025-04-10 11:17:30"
$data[0,6] = "2025-04-10 11:17:30"
$data[0,7] = "2025-04-10 12:27:39"
$data[0,8] = 3859
$data[0,9] = "bobina"
$data[0,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$data[0,11] = 8
$data[0,12] = 70
$data[0,13] = 4

$data[1,0] = 251109
$data[1,1] = "R6"
$data[1,2] = 112
$data[1,3] = 266.5915492957747
$data[1,4] = "2025-04-10 13:25:00"
$data[1,5] = "2025-04-11 07:17:00"
$data[1,6] = "2025-04-11 07:17:00"
$data[1,7] = "2025-04-11 11:43:35"
$data[1,8] = 18928
$data[1,9] = "bobina"
$data[1,10] = "R6"
$data[1,11] = 16
$data[1,12] = 70
$data[1,13] = 4

$data[2,0] = 251204
$data[2,1] = "R9"
$data[2,2] = 30
$data[2,3] = 111.0579710144928
$data[2,4] = "2025-04-10 07:22:00"
$data[2,5] = "2025-04-10 07:52:00"
$data[2,6] = "2025-04-10 07:52:00"
$data[2,7] = "2025-04-10 09:43:03"
$data[2,8] = 7663
$data[2,9] = "bobina"
$data[2,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$data[2,11] = 2
$data[2,12] = 70
$data[2,13] = 4

$data[3,0] = 251088
$data[3,1] = "R9"
$data[3,2] = 30
$data[3,3] = 89.6376811594203
$data[3,4] = "2025-04-10 09:43:03"
$data[3,5] = "2025-04-10 10:13:03"
$data[3,6] = "2025-04-10 10:13:03"
$data[3,7] = "2025-04-10 11:42:41"
$data[3,8] = 6185
$data[3,9] = "bobina"
$data[3,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$data[3,11] = 3
$data[3,12] = 70
$data[3,13] = 4

$data[4,0] = 251550
$data[4,1] = "R9"
$data[4,2] = 25
$data[4,3] = 516.6811594202899
$data[4,4] = "2025-04-10 11:42:41"
$data[4,5] = "2025-04-10 12:07:41"
$data[4,6] = "2025-04-10 12:07:41"
$data[4,7] = "2025-04-11 12:44:22"
$data[4,8] = 35651
$data[4,9] = "bobina"
$data[4,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$data[4,11] = 3
$data[4,12] = 70
$data[4,13] = 4

$data[5,0] = 235572
$data[5,1] = "R9"
$data[5,2] = 35
$data[5,3] = 144.3188405797102
$data[5,4] = "2025-04-11 12:44:22"
$data[5,5] = "2025-04-11 13:19:22"
$data[5,6] = "2025-04-11 13:19:22"
$data[5,7] = "2025-04-14 07:43:41"
$data[5,8] = 9958
$data[5,9] = "bobina"
$data[5,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R6 ;R9"
$data[5,11] = 5
$data[5,12] = 70
$data[5,13] = 4

$data[6,0] = 250866
$data[6,1] = "R9"
$data[6,2] = 30
$data[6,3] = 74.14492753623189
$data[6,4] = "2025-04-14 07:43:41"
$data[6,5] = "2025-04-14 08:13:41"
$data[6,6] = "2025-04-14 08:13:41"
$data[6,7] = "2025-04-14 09:27:50"
$data[6,8] = 5116
$data[6,9] = "bobina"
$data[6,10] = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$data[6,11] = 6
$data[6,12] = 70
$data[6,13] = 4

$ws.Range("A3:N9").Value = $data

# Remove the rows that are no longer needed; the sheet now ends at row 9.
$ws.Rows("10:14").Delete()
